# Applies the Balance.xlsx data refresh:
#  - swap the "Max aft Xcg MAC" / "Max forward Xcg MAC" labels on GLOBAL RESULTS
#  - update recomputed numeric results on GLOBAL RESULTS, WING, FUEL TANK and
#    LANDING GEARS sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# GLOBAL RESULTS
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws.Range("A22").Value = "Max forward Xcg MAC"
$ws.Range("A23").Value = "Max aft Xcg MAC"

$ws.Range("C2").Value  = -16.329351801878612
$ws.Range("C3").Value  = 19.75693186198474
$ws.Range("C4").Value  = -0.5174378876703667
$ws.Range("C6").Value  = 88.1977944981336
$ws.Range("C7").Value  = 22.837106327893174
$ws.Range("C8").Value  = -0.0282959666618242
$ws.Range("C10").Value = 88.1977944981336
$ws.Range("C11").Value = 22.837106327893174
$ws.Range("C12").Value = -0.0282959666618242
$ws.Range("C14").Value = 3.4924007008475004
$ws.Range("C15").Value = 20.341033294862434
$ws.Range("C16").Value = -0.01896605115289715
$ws.Range("C18").Value = 12.098162583058741
$ws.Range("C19").Value = 20.594625294142645
$ws.Range("C20").Value = -0.3240984506076898
$ws.Range("C22").Value = 0.25831066811662906
$ws.Range("C23").Value = 90.97143260202606

# ---------------------------------------------------------------------------
# WING
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")

$ws.Range("C2").Value  = 2.0862310247028493
$ws.Range("C3").Value  = 5.766295695644638
$ws.Range("C6").Value  = 21.176231024702844
$ws.Range("C7").Value  = 5.766295695644637
$ws.Range("C11").Value = 2.0862310247028493
$ws.Range("C14").Value = 5.766295695644638

# ---------------------------------------------------------------------------
# FUEL TANK
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FUEL TANK")

$ws.Range("C2").Value = 2.332307524701931
$ws.Range("C6").Value = 21.422307524701928

# ---------------------------------------------------------------------------
# LANDING GEARS
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LANDING GEARS")

$ws.Range("C2").Value = 18.071707001519314
